{"js": "// Applies the CV content updates described by the diff:\n//  1. Trim the \"direct $XXM+ portfolio with\" clause from the D&A bullet.\n//  2. Remove the \"Designed credit risk AI models...\" bullet entirely.\n//  3. Rename \"Various Companies\" to \"Microsoft, UTU & Others\".\n//  4. Rename the role title from \"...Technical Consulting\" to \"...Technical Leadership\".\n//  5. Replace the \"Progressive roles...\" summary line with the new one.\n//  6. Drop the trailing \"JAPAC Hackathon Winner.\" clause from the CatchMe bullet.\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${findText}\", found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Shorten the CSAT bullet.\nawait replaceOnce(\n  \"Built D&A practice from 0 to 1 across 6 countries; direct $XXM+ portfolio with 97% CSAT\",\n  \"Built D&A practice from 0 to 1 across 6 countries; 97% CSAT\"\n);\n\n// 2. Delete the \"Designed credit risk AI models...\" bullet paragraph.\nconst creditResults = body.search(\n  \"Designed credit risk AI models improving accuracy by 15% with alternative data sources\",\n  { matchCase: true }\n);\ncreditResults.load(\"items\");\nawait context.sync();\nif (creditResults.items.length !== 1) {\n  throw new Error(\n    `Expected exactly 1 match for the credit-risk bullet, found ${creditResults.items.length}`\n  );\n}\ncreditResults.items[0].paragraphs.load(\"items\");\nawait context.sync();\ncreditResults.items[0].paragraphs.items[0].delete();\nawait context.sync();\n\n// 3. Rename the \"Various Companies\" heading.\nawait replaceOnce(\"Various Companies\", \"Microsoft, UTU & Others\");\n\n// 4. Rename the role title.\nawait replaceOnce(\n  \"Software Engineering & Technical Consulting\",\n  \"Software Engineering & Technical Leadership\"\n);\n\n// 5. Replace the role summary sentence.\nawait replaceOnce(\n  \"Progressive roles in software development, systems integration, and consulting in financial services and algorithmic trading.\",\n  \"Windows Kernel development (Microsoft), payment systems (UTU Singapore), founded Truckaurbus B2B marketplace.\"\n);\n\n// 6. Drop the trailing \"JAPAC Hackathon Winner.\" clause.\nawait replaceOnce(\n  \"Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency. JAPAC Hackathon Winner.\",\n  \"Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency.\"\n);\n", "ps1": "# Applies the CV content updates described by the diff:\n#  1. Trim the \"direct $XXM+ portfolio with\" clause from the D&A bullet.\n#  2. Remove the \"Designed credit risk AI models...\" bullet entirely.\n#  3. Rename \"Various Companies\" to \"Microsoft, UTU & Others\".\n#  4. Rename the role title from \"...Technical Consulting\" to \"...Technical Leadership\".\n#  5. Replace the \"Progressive roles...\" summary line with the new one.\n#  6. Drop the trailing \"JAPAC Hackathon Winner.\" clause from the CatchMe bullet.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Find.Execute could not locate text: $findText\"\n    }\n}\n\n# 1. Shorten the CSAT bullet.\nReplace-Once \"Built D&A practice from 0 to 1 across 6 countries; direct `$XXM+ portfolio with 97% CSAT\" \"Built D&A practice from 0 to 1 across 6 countries; 97% CSAT\"\n\n# 2. Delete the \"Designed credit risk AI models...\" bullet paragraph.\n$deleted = $false\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*Designed credit risk AI models improving accuracy by 15% with alternative data sources*\") {\n        $p.Range.Delete()\n        $deleted = $true\n        break\n    }\n}\nif (-not $deleted) {\n    throw \"Could not find the credit-risk bullet paragraph to delete\"\n}\n\n# 3. Rename the \"Various Companies\" heading.\nReplace-Once \"Various Companies\" \"Microsoft, UTU & Others\"\n\n# 4. Rename the role title.\nReplace-Once \"Software Engineering & Technical Consulting\" \"Software Engineering & Technical Leadership\"\n\n# 5. Replace the role summary sentence.\nReplace-Once \"Progressive roles in software development, systems integration, and consulting in financial services and algorithmic trading.\" \"Windows Kernel development (Microsoft), payment systems (UTU Singapore), founded Truckaurbus B2B marketplace.\"\n\n# 6. Drop the trailing \"JAPAC Hackathon Winner.\" clause.\nReplace-Once \"Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency. JAPAC Hackathon Winner.\" \"Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency.\"\n\nWrite-Output \"done\"\n"}
